# Regenerate save_data to use K (strikeouts) instead of Strike# for column G,
# writing the newly calculated s_vals (K) for each outing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G)
$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 3
    7  = 1
    8  = 2
    9  = 3
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    20 = 1
    21 = 3
    22 = 1
    23 = 1
    24 = 3
    25 = 1
    26 = 0
    27 = 3
    28 = 1
    29 = 0
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 2
    35 = 0
    36 = 1
    37 = 1
    38 = 2
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    45 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
